$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$ws.Cells.Item(1,1).Value  = "code"
$ws.Cells.Item(1,2).Value  = "name"
$ws.Cells.Item(1,3).Value  = "descr"
$ws.Cells.Item(1,4).Value  = "lang_code"
$ws.Cells.Item(1,5).Value  = "is_active"
$ws.Cells.Item(1,6).Value  = "cr_by"
$ws.Cells.Item(1,7).Value  = "cr_dtimes"
$ws.Cells.Item(1,8).Value  = "upd_by"
$ws.Cells.Item(1,9).Value  = "upd_dtimes"
$ws.Cells.Item(1,10).Value = "is_deleted"
$ws.Cells.Item(1,11).Value = "del_dtimes"

# ---- Row 2 : REG / Permanent ----
$ws.Cells.Item(2,1).Value  = "REG"
$ws.Cells.Item(2,2).Value  = "Permanent"
$ws.Cells.Item(2,3).Value  = "Centre d'Enregistrement Permanent"
$ws.Cells.Item(2,4).Value  = "fra"
$ws.Cells.Item(2,5).Value  = $true
$ws.Cells.Item(2,6).Value  = "superadmin"
$ws.Cells.Item(2,7).Value  = 45079.57731385417
$ws.Cells.Item(2,7).NumberFormat = "mm:ss.0"
$ws.Cells.Item(2,8).Value  = "NULL"
$ws.Cells.Item(2,9).Value  = "NULL"
$ws.Cells.Item(2,10).Value = $false
$ws.Cells.Item(2,11).Value = "NULL"

# ---- Row 3 : CEP / Provisoire ----
$ws.Cells.Item(3,1).Value  = "CEP"
$ws.Cells.Item(3,2).Value  = "Provisoire"
$ws.Cells.Item(3,3).Value  = "Centre d'Enregistrement Provisoire"
$ws.Cells.Item(3,4).Value  = "fra"
$ws.Cells.Item(3,5).Value  = $true
$ws.Cells.Item(3,6).Value  = "superadmin"
$ws.Cells.Item(3,7).Value  = 45079.57731385417
$ws.Cells.Item(3,7).NumberFormat = "mm:ss.0"
$ws.Cells.Item(3,8).Value  = "NULL"
$ws.Cells.Item(3,9).Value  = "NULL"
$ws.Cells.Item(3,10).Value = $false
$ws.Cells.Item(3,11).Value = "NULL"

# ---- Row 4 : CEI / Itinérant (mojibake'd name/descr, matching source data) ----
$ws.Cells.Item(4,1).Value  = "CEI"
$ws.Cells.Item(4,2).Value  = "ItinÃ©rant"
$ws.Cells.Item(4,3).Value  = "Centre d'Enregistrement ItinÃ©rant"
$ws.Cells.Item(4,4).Value  = "fra"
$ws.Cells.Item(4,5).Value  = $true
$ws.Cells.Item(4,6).Value  = "superadmin"
$ws.Cells.Item(4,7).Value  = 45079.57731385417
$ws.Cells.Item(4,7).NumberFormat = "mm:ss.0"
$ws.Cells.Item(4,8).Value  = "NULL"
$ws.Cells.Item(4,9).Value  = "NULL"
$ws.Cells.Item(4,10).Value = $false
$ws.Cells.Item(4,11).Value = "NULL"

# ---- Selection matches the saved view state ----
[void]$ws.Range("C13").Select()
